$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main Pages")
$ws.Range("A1").Value = "test"
